$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.130.95'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.35%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.294.55'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.28%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.63'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.23%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.66'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.36%  '

# Row 7
$ws.Range("E7").Value = '  +0.13%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.290.98'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.34%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.515'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.62%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.148'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.91%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.49'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.91%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.464'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.08%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000244'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.00%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.37'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.41%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.851.66'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.80%  '

# Row 16
$ws.Range("E16").Value = '  +0.93%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.305.46'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.76%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.343.83'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.12%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.79'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.64%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '475.65'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.31%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.88'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.06%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.730'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.17%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.91'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.33%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.79'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.99%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.93'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.21%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.76'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.77%  '

# Row 28
$ws.Range("E28").Value = '  -0.02%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.12'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.40%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.06'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.11%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.12'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.09%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.40'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.91%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.104'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.30%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.49'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.16%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.08'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.05%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.98'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.08%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.18'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.94%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0729'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.88%  '

# Row 39
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.127.10'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.83%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0398'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.47%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '425.06'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.30%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.118'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +7.83%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.25'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.20%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.67'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.89%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.261'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.39%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.19'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.36%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '36.29'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +8.92%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.85'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.76%  '

# Row 49
$ws.Range("B49").Value = 'USDe'
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.999'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.07%  '

# Row 50
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.06'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.16%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.29'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.82%  '
